$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link text cells (rows 36-39 reordering)
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

# Force Price/Volume columns to text so numeric-looking strings are not
# auto-converted to numbers by Excel, matching the original inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("D2").Value = '98.923.21'
$ws.Range("E2").Value = '  +1.80%  '
$ws.Range("D3").Value = '3.313.06'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '255.77'
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").Value = '625.30'
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("E7").Value = '  +32.67%  '
$ws.Range("D8").Value = '0.409'
$ws.Range("E8").Value = '  +6.76%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("D10").Value = '0.970'
$ws.Range("E10").Value = '  +23.08%  '
$ws.Range("D11").Value = '3.311.11'
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("E13").Value = '  +12.19%  '
$ws.Range("D14").Value = '98.620.27'
$ws.Range("E14").Value = '  +1.64%  '
$ws.Range("D15").Value = '0.0000251'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '3.930.80'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '5.48'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '3.305.49'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").Value = '3.48'
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("D20").Value = '15.75'
$ws.Range("E20").Value = '  +5.73%  '
$ws.Range("D21").Value = '6.30'
$ws.Range("E21").Value = '  +8.92%  '
$ws.Range("D22").Value = '488.55'
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("D23").Value = '9.50'
$ws.Range("E23").Value = '  +3.36%  '
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("D25").Value = '5.63'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '88.80'
$ws.Range("E26").Value = '  +1.38%  '
$ws.Range("D27").Value = '12.08'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = '0.302'
$ws.Range("E28").Value = '  +27.22%  '
$ws.Range("D29").Value = '3.482.58'
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '0.136'
$ws.Range("E31").Value = '  +12.49%  '
$ws.Range("D32").Value = '0.189'
$ws.Range("E32").Value = '  +3.75%  '
$ws.Range("E33").Value = '  +9.06%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = '28.01'
$ws.Range("E35").Value = '  +2.62%  '
$ws.Range("D36").Value = '0.471'
$ws.Range("E36").Value = '  +5.25%  '
$ws.Range("D37").Value = '7.21'
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("D38").Value = '0.148'
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = '1.95'
$ws.Range("E39").Value = '  +0.86%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '490.44'
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("E42").Value = '  +3.77%  '
$ws.Range("D43").Value = '1.24'
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").Value = '0.790'
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D46").Value = '3.12'
$ws.Range("E46").Value = '  -5.48%  '
$ws.Range("D47").Value = '159.40'
$ws.Range("E47").Value = '  -1.29%  '
$ws.Range("D48").Value = '1.96'
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("D49").Value = '7.32'
$ws.Range("E49").Value = '  +16.19%  '
$ws.Range("D50").Value = '0.846'
$ws.Range("E50").Value = '  +6.13%  '
$ws.Range("D51").Value = '4.71'
$ws.Range("E51").Value = '  +5.02%  '

# Restore default (unstyled) cell formatting so the style index matches the original.
$ws.Range("D2:E51").ClearFormats()

